$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: 7569 -> 7295
$ws.Range("C2").Value = 7295

# C3:C252: 7569 -> 7293
$ws.Range("C3:C252").Value = 7293
